$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '43.779.52'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.283.93'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.96%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '120.19'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '267.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("E7").Value = '  +4.57%  '
$ws.Range("E8").Value = '  +0.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.629'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.65%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.57'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.94%  '
$ws.Range("E11").Value = '  +2.99%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.28'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +6.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.65'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.25%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.921'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.08%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.626.74'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.285.15'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '43.808.41'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.67%  '
$ws.Range("E19").Value = '  +3.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.96'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.86%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +2.06%  '
$ws.Range("E22").Value = '  +1.08%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.56'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.82%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.65'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.55%  '
$ws.Range("E25").Value = '  +2.29%  '
$ws.Range("E26").Value = '  +6.38%  '
$ws.Range("E27").Value = '  +1.84%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '42.73'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.83%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("E30").Value = '  +0.26%  '
$ws.Range("E31").Value = '  +0.74%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.69'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.04%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0933'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.81'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.61%  '
$ws.Range("E35").Value = '  +4.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.28'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +14.36%  '
$ws.Range("E37").Value = '  +10.92%  '
$ws.Range("E38").Value = '  +1.15%  '
$ws.Range("E39").Value = '  +3.69%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.57'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.44'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.15%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.79'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.46%  '
$ws.Range("E43").Value = '  +2.08%  '
$ws.Range("E44").Value = '  +0.04%  '
$ws.Range("E45").Value = '  +2.23%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '5.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.40%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '76.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +48.91%  '
$ws.Range("E48").Value = '  +3.15%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.57'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.23%  '
$ws.Range("E50").Value = '  +1.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '102.57'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.84%  '
